# TC05_Bento_Filter_Diagnosis-InfilDuctalCarcinoma.xlsx
# "Added Samples and Files Tab to all tests"
#
# Adds two new rows (SamplesTab, FilesTab) to the "startup" sheet, alongside
# the existing CasesTab row, each carrying its own Neo4j/Cypher query (col B),
# a shared StatQuery (col C), and the shared Neo4jData/WebData file names
# (cols D/E). Also makes a small textual fix to the existing CasesTab query
# (col B) and refreshes column widths / selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Updated CasesTab query (row 2, column B) - minor fix to the WHERE clause
# (dropped a stray extra "]" and a trailing space before the newline).
# ---------------------------------------------------------------------
$casesQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
 WHERE ss.disease_subtype IN ["Infiltrating Ductal Carcinoma"] 
return ss.study_subject_id as `Case ID`,
       p.program_acronym as `Program Code`,
        p.program_id as Program_ID,
       s.study_acronym as `Arm`,
       ss.disease_subtype as `Diagnosis`,
       sf.grouped_recurrence_score AS `Recurrence Score`,
       d.tumor_size_group AS `tumor_size`,
       d.er_status AS `ER Status`,
       d.pr_status AS `PR Status`,
       demo.age_at_index AS `Age (years)`,
demo.survival_time AS `Survival (days)`
'@

# ---------------------------------------------------------------------
# New SamplesTab query (row 3, column B)
# ---------------------------------------------------------------------
$samplesQuery = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE ss.disease_subtype IN ["Infiltrating Ductal Carcinoma"] 
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`
'@

# ---------------------------------------------------------------------
# New FilesTab query (row 4, column B)
# ---------------------------------------------------------------------
$filesQuery = @'
MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE ss.disease_subtype IN ["Infiltrating Ductal Carcinoma"] 
WITH
        f, parent,p, ss, d,tp, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`
    order by f.file_name
'@

# ---------------------------------------------------------------------
# New shared StatQuery (column C) used by all three rows
# ---------------------------------------------------------------------
$statQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE ss.disease_subtype IN ["Infiltrating Ductal Carcinoma"]
WITH ss
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
'@

$neo4jFile = 'TC05_Bento_Filter_Diagnosis-InfilDuctalCarcinoma_Neo4jData.xlsx'
$webFile   = 'TC05_Bento_Filter_Diagnosis-InfilDuctalCarcinoma_WebData.xlsx'

# ---------------------------------------------------------------------
# Row 2 : CasesTab (existing row, query text refreshed)
# ---------------------------------------------------------------------
$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("D2").Value = $neo4jFile
$ws.Range("E2").Value = $webFile
$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true

# ---------------------------------------------------------------------
# Row 3 : SamplesTab (new row)
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = $neo4jFile
$ws.Range("E3").Value = $webFile
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# ---------------------------------------------------------------------
# Row 4 : FilesTab (new row)
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "FilesTab"
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery
$ws.Range("D4").Value = $neo4jFile
$ws.Range("E4").Value = $webFile
$ws.Range("B4").WrapText = $true
$ws.Range("C4").WrapText = $true

# ---------------------------------------------------------------------
# Row heights (best-effort match to the saved workbook's wrapped heights)
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 316.8
$ws.Rows.Item(3).RowHeight = 345.6
$ws.Rows.Item(4).RowHeight = 409.6

# ---------------------------------------------------------------------
# Column widths (best-effort match; engine rounds to its own granularity)
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 11.94
$ws.Columns.Item(2).ColumnWidth = 75.28
$ws.Columns.Item(3).ColumnWidth = 58.61
$ws.Columns.Item(4).ColumnWidth = 64.72
$ws.Columns.Item(5).ColumnWidth = 63.28

# ---------------------------------------------------------------------
# Selection, matching the saved workbook (B4 selected)
# ---------------------------------------------------------------------
$null = $ws.Range("B4").Select()
